# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for rows 2-18 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 2
    9  = 1
    10 = 2
    11 = 0
    12 = 1
    13 = 0
    14 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
